# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# worksheets. Values below are the latest computed figures; some rows also
# gain or lose an M (LeveProfitNQ) / N (LeveProfitHQ) cell depending on
# whether that computation is meaningful for the row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10468

$ws.Range("H28").Value = 922.0909
$ws.Range("I28").Value = 551.25
$ws.Range("J28").Value = 1911
$ws.Range("K28").Value = 551.25
$ws.Range("L28").Value = 1911
$ws.Range("M28").Value = -66.25
$ws.Range("N28").Value = -2881

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H46").Value = 1081.9
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1091
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 3273
$ws.Range("M46").Value = -2881
$ws.Range("N46").Value = -3511

$ws.Range("H57").Value = 23870
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 23870
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 71610
$ws.Range("N57").Value = -72608

$ws.Range("H60").Value = 1081.9
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 1091
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 3273
$ws.Range("M60").Value = -2516
$ws.Range("N60").Value = -4241

$ws.Range("H80").Value = 254.27272
$ws.Range("I80").Value = 236.19048
$ws.Range("J80").Value = 285.91666
$ws.Range("K80").Value = 708.5714400000001
$ws.Range("L80").Value = 857.7499799999999
$ws.Range("M80").Value = 289.4285599999999
$ws.Range("N80").Value = -2853.74998

$ws.Range("H83").Value = 254.27272
$ws.Range("I83").Value = 236.19048
$ws.Range("J83").Value = 285.91666
$ws.Range("K83").Value = 2125.71432
$ws.Range("L83").Value = 2573.24994
$ws.Range("M83").Value = 2866.28568
$ws.Range("N83").Value = -12557.24994

$ws.Range("H137").Value = 1951
$ws.Range("I137").Value = 867.3333
$ws.Range("J137").Value = 2601.2
$ws.Range("K137").Value = 2601.9999
$ws.Range("L137").Value = 7803.599999999999
$ws.Range("M137").Value = -51.9998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 6504.5
$ws.Range("I30").Value = 4009
$ws.Range("J30").Value = 9000
$ws.Range("K30").Value = 4009
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = -3859

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20010
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 20010
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 20010
$ws.Range("N19").Value = -20356

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H55").Value = 19332.666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 19332.666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 19332.666
$ws.Range("N55").Value = -19878.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2570
$ws.Range("I19").Value = 2084
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 2084
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -1914
$ws.Range("N19").Value = -5340

$ws.Range("H24").Value = 2570
$ws.Range("I24").Value = 2084
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 2084
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -1914
$ws.Range("N24").Value = -5340

$ws.Range("H31").Value = 1196.7742
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 3875
$ws.Range("K31").Value = 800
$ws.Range("L31").Value = 3875
$ws.Range("M31").Value = -505
$ws.Range("N31").Value = -4465

$ws.Range("H34").Value = 1196.7742
$ws.Range("I34").Value = 800
$ws.Range("J34").Value = 3875
$ws.Range("K34").Value = 800
$ws.Range("L34").Value = 3875
$ws.Range("M34").Value = -598
$ws.Range("N34").Value = -4279

$ws.Range("H81").Value = 37080
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 37080
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 37080
$ws.Range("N81").Value = -39076

$ws.Range("H84").Value = 37080
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 37080
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 111240
$ws.Range("N84").Value = -121224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2757
$ws.Range("I87").Value = 2757
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 8271
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -7023

$ws.Range("H90").Value = 2757
$ws.Range("I90").Value = 2757
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 24813
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -18573

$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 7500
$ws.Range("N100").Value = -9122

$ws.Range("H114").Value = 351.94116
$ws.Range("I114").Value = 215.41667
$ws.Range("J114").Value = 679.6
$ws.Range("K114").Value = 646.25001
$ws.Range("L114").Value = 2038.8
$ws.Range("M114").Value = 2607.74999
$ws.Range("N114").Value = -8546.799999999999

$ws.Range("H120").Value = 5493
$ws.Range("I120").Value = 3239.5
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 9718.5
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -4880.5

$ws.Range("H129").Value = 1338.7646
$ws.Range("I129").Value = 895
$ws.Range("J129").Value = 1733.2222
$ws.Range("K129").Value = 2685
$ws.Range("L129").Value = 5199.6666
$ws.Range("M129").Value = 2315

$ws.Range("H131").Value = 899.4400000000001
$ws.Range("I131").Value = 770
$ws.Range("J131").Value = 904.8333
$ws.Range("K131").Value = 2310
$ws.Range("L131").Value = 2714.4999
$ws.Range("M131").Value = 2730
$ws.Range("N131").Value = -12794.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 25913.334
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 25913.334
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 25913.334
$ws.Range("N21").Value = -26259.334

$ws.Range("H30").Value = 25913.334
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 25913.334
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 25913.334
$ws.Range("N30").Value = -26123.334

$ws.Range("H31").Value = 665
$ws.Range("I31").Value = 665
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 665
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -373

$ws.Range("H37").Value = 665
$ws.Range("I37").Value = 665
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 665
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -388

$ws.Range("H41").Value = 418.25
$ws.Range("I41").Value = 418.25
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 418.25
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -63.25
$ws.Range("N41").ClearContents()

$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -4417

$ws.Range("H110").Value = 30702
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 30702
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 30702
$ws.Range("N110").Value = -38882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3000
$ws.Range("I4").Value = 3000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2887
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 3000
$ws.Range("I28").Value = 3000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2768
$ws.Range("N28").ClearContents()

$ws.Range("H30").Value = 1450
$ws.Range("I30").Value = 933.3333
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 933.3333
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = -825.3333
$ws.Range("N30").Value = -3216

$ws.Range("H35").Value = 17958.25
$ws.Range("I35").Value = 7277.6665
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 7277.6665
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -6941.6665

$ws.Range("H37").Value = 3000
$ws.Range("I37").Value = 3000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2893
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H46").Value = 44976.668
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 44976.668
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 44976.668
$ws.Range("N46").Value = -45438.668

$ws.Range("H134").Value = 44976.668
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 44976.668
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 134930.004
$ws.Range("N134").Value = -140000.004
